$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new test-case row (row 50) below the existing last row (49) ---
# Copy the whole row's formatting/values down first so the new row inherits
# the same cell styles (borders / vertical alignment / wrap) as its template.
$ws.Range("A49:L49").Copy($ws.Range("A50:L50"))

# Column K holds a value on the new row (unlike row 49, where K is blank), so
# pull its formatting from a neighbouring "has value" cell (J49) instead.
$ws.Range("J49").Copy($ws.Range("K50"))

# Fill in the new row's content. The order below matches the order the new
# strings were appended to the shared-string table in the source edit:
#   ?name=wang&...                                            -> QUERYSTRING (G)
#   status=200||hits.primaryName=wang||...                    -> VALIDATIONS (J)
#   hits[0].publicationYearRangeMin||hits.publicationYearRangeMax -> STORE (K)
#   Verify that user should be able to filter ... given an order -> DESCRIPTION (B)
#   WAT-413                                                    -> TESTNAME (A)
$ws.Range("G50").Value = "?name=wang&affiliation=china&filter=name&category=physics&filter=affiliation&filter=catagory&sort=year&order=asc&limit=10"
$ws.Range("J50").Value = "status=200||hits.primaryName=wang||filters.category=physics||hits.affiliation=china"
$ws.Range("K50").Value = "hits[0].publicationYearRangeMin||hits.publicationYearRangeMax"
$ws.Range("B50").Value = "Verify that user should be able to filter values for a set of authors  provided with all the mandatory inputs along with sorting given an order"
$ws.Range("A50").Value = "WAT-413"

# The rest of the row mirrors row 49 (same HOST / APIPATH / METHOD), already
# brought over by the row copy above: C50=WOSAUTHORRECOMMEND, D50=/author/search,
# E50=GET. F50, H50 and I50 stay blank, same as row 49.

# New row is taller (wraps to 3 lines) than the template row.
$ws.Rows(50).RowHeight = 45

# Move the selection the way the author's session ended up after adding the
# row (scrolled down, clicked into the empty row below the new data).
$ws.Range("A53").Select() | Out-Null
